$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This reverts the "Done with phase 1" commit: the Bill-of-materials sheet
# goes back to its earlier (phase 0) content - different part numbers,
# quantities, descriptions, "where to get" notes, and no unit prices filled
# in (so the Price column mostly computes to 0, except a couple of rows).
# ---------------------------------------------------------------------------

# Clear out the whole data area first so no stale values/formulas survive
# in cells that should end up blank.
$ws.Range("A2:G11").ClearContents()

# The unit-price column (E) is blank for every row except row 4 in the
# reverted version - fully clear its formatting too (not just the value)
# so those cells disappear from the sheet entirely instead of lingering
# as empty-but-currency-formatted cells. Row 4 keeps its existing
# currency-formatted cell (its value is overwritten further down), so
# leave it alone here.
$ws.Range("E2:E3").Clear()
$ws.Range("E5:E11").Clear()

# --- Row 1 (headers) -------------------------------------------------------
$ws.Range("A1").Value = "Part No"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Qty"
$ws.Range("D1").Value = "Use"
$ws.Range("E1").Value = "unit price"
$ws.Range("F1").Value = "Price"
$ws.Range("G1").Value = "Where to get"

# --- Row 2 -------------------------------------------------------------
$ws.Range("B2").Value = "5VDC motor"
$ws.Range("C2").Value = 1
$ws.Range("F2").Formula = "=E2*C2"
$ws.Range("G2").Value = "have it"

# --- Row 3 -------------------------------------------------------------
$ws.Range("A3").Value = "PIC18F452"
$ws.Range("B3").Value = "40-pin PIC microcontroller"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "nodes/controller"
$ws.Range("F3").Formula = "=E3*C3"
$ws.Range("G3").Value = "have it"

# --- Row 4 -------------------------------------------------------------
$ws.Range("A4").Value = "OSC-20MHZ"
$ws.Range("B4").Value = "Xtal Oscillator"
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = 2.4
$ws.Range("F4").Formula = "=E4*C4"
$ws.Range("G4").Value = "EE store"

# --- Row 5 -------------------------------------------------------------
$ws.Range("A5").Value = "TIP22"
$ws.Range("B5").Value = "Mosfet"
$ws.Range("D5").Value = "Motor driver"
$ws.Range("F5").Formula = "=E5*C5"
$ws.Range("G5").Value = "haveit ?"

# --- Row 6 -------------------------------------------------------------
$ws.Range("B6").Value = "GALs"
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = "I/O buffer"
$ws.Range("F6").Formula = "=E6*C6"
$ws.Range("G6").Value = "have it"

# --- Row 7 -------------------------------------------------------------
$ws.Range("B7").Value = "Lots of Wires"
$ws.Range("C7").Value = "infinity"
$ws.Range("F7").Formula = "=E7*C7"

# --- Row 8 -------------------------------------------------------------
$ws.Range("A8").Value = "MAX232"
$ws.Range("B8").Value = "RS-232 level-shifter"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "UART-rs232 level shifter"
$ws.Range("F8").Value = 0.74
$ws.Range("G8").Value = "EE store"

# --- Row 9 -------------------------------------------------------------
$ws.Range("A9").Value = "CY7C128A"
$ws.Range("B9").Value = "SRAM"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "Store measurment data"
$ws.Range("F9").Formula = "=E9*C9"
$ws.Range("G9").Value = "have them"

# --- Row 10 ------------------------------------------------------------
$ws.Range("B10").Value = "Resistors"
$ws.Range("F10").Formula = "=E10*C10"

# --- Row 11 ------------------------------------------------------------
$ws.Range("B11").Value = "Capacitors"
$ws.Range("F11").Formula = "=E11*C11"

# --- Column B width reverts to its earlier (narrower) best-fit size -------
$ws.Columns.Item(2).ColumnWidth = 24.43

# --- Selection reverts to its earlier position -----------------------------
$ws.Range("C14").Select()

$wb.Application.Calculate()
